$wb = $excel.ActiveWorkbook

# Insert a new "PicklistValues" worksheet right after "Users"
$usersSheet = $wb.Worksheets.Item("Users")
$picklistSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $usersSheet)
$picklistSheet.Name = "PicklistValues"

# Populate the picklist values (entered in this order so the shared-string
# table ends up in the same order as the source workbook)
$picklistSheet.Range("A1").Value = "EngPotentialRoundTrip PicklistValues"
$picklistSheet.Range("A4").Value = "Neither subject nor buyer are round trip"
$picklistSheet.Range("A2").Value = "Subject is a potential round trip"
$picklistSheet.Range("A3").Value = "Buyer is a potential round trip"

# Header formatting + column width
$picklistSheet.Range("A1").Font.Bold = $true
$picklistSheet.Columns.Item(1).ColumnWidth = 35.21875

# Update selections to match the saved state
$usersSheet.Range("G15").Select() | Out-Null
$picklistSheet.Range("A8").Select() | Out-Null
$picklistSheet.Activate() | Out-Null
